$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns at J:K - this shifts the existing
#    "Marca do celular".."Ativo Notebook" block (J:Q) over to L:S
$ws.Columns("J:K").Insert()

# 2. Add the three new trailing columns (Notebook component info)
#    Fill header + data in the same order the original author appears to
#    have used, row by row, so that the resulting workbook matches the
#    expected content.
$ws.Range("T1").Value = "Armazenamento Notebook"
$ws.Range("U1").Value = "Processador Notebook"
$ws.Range("V1").Value = "Memoria RAM Notebook"

# Armazenamento Notebook (storage) column, filled top to bottom
$ws.Range("T2").Value = "HDD 500GB"
$ws.Range("T3").Value = "HDD 500GB"
$ws.Range("T4").Value = "SSD 256GB"
$ws.Range("T7").Value = "SSD 256GB"
$ws.Range("T8").Value = "SSD 500GB"
$ws.Range("T10").Value = "SSD 256GB"

# Processador Notebook (CPU) column, filled top to bottom
$ws.Range("U2").Value = "i5"
$ws.Range("U3").Value = "i5"
$ws.Range("U4").Value = "i5"
$ws.Range("U7").Value = "i7"
$ws.Range("U8").Value = "i7"
$ws.Range("U10").Value = "i5"

# Memoria RAM Notebook (RAM) column
$ws.Range("V4").Value = "16GB"
$ws.Range("V7").Value = "16GB"
$ws.Range("V8").Value = "16GB"
$ws.Range("V10").Value = "16GB"
$ws.Range("V2").Value = "8GB"
$ws.Range("V3").Value = "8GB"

# 3. Add the new "Plano" / "Status" columns in the gap created above
$ws.Range("J1").Value = "Plano"
$ws.Range("K1").Value = "Status"

$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3

$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 3

$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3

$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 5

$ws.Range("J11").Value = 3
$ws.Range("K11").Value = 5

# 4. Match the header styling (grey fill) used by the rest of row 1
$ws.Range("A1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)
$ws.Range("T1:V1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 5. Sheet view tweaks
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("J36").Select()

Write-Host "done"
